$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 112043839
$ws.Range("B2").Value = 78713
$ws.Range("Q2").Value = 547969
$ws.Range("R2").Value = 6960405

# --- Row 3 (becomes the "Tretåig hackspett" observation) ---
$ws.Range("A3").Value = 112043807
$ws.Range("B3").Value = 56430
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
# A lone apostrophe yields a "quote-prefixed" empty text cell - the closest
# COM-surface equivalent of the source file's empty <is><t/></is> cell.
$ws.Range("K3").Value = "'"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'"
$ws.Range("N3").Value = "'"
$ws.Range("Q3").Value = 547961
$ws.Range("R3").Value = 6960421
$ws.Range("AC3").Value = "ringhack"

# --- Row 4 (becomes the "Lunglav" observation) ---
$ws.Range("A4").Value = 112043819
$ws.Range("B4").Value = 78713
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("Q4").Value = 547979
$ws.Range("R4").Value = 6960195
$ws.Range("AC4").ClearContents()
